# Grouptime.xlsx final update: split the single "Sheet1" into three tabs
# ("Group 1", "Group 2", "Other") and leave a couple of cursor/selection
# breadcrumbs behind, matching the author's last save before moving on to
# another project.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ------------------------------------------------
$group1 = $wb.Worksheets.Item(1)
$group1.Name = "Group 1"

# --- Add "Group 2" right after "Group 1" --------------------------------------
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "Group 2"
$group1 = $wb.Worksheets.Item("Group 1")
$group2 = $wb.Worksheets.Item("Group 2")
$group2.Move([Type]::Missing, $group1)

# --- Add "Other" right after "Group 2" ----------------------------------------
$null = $wb.Worksheets.Add()
$wb.Worksheets.Item(1).Name = "Other"
$group2 = $wb.Worksheets.Item("Group 2")
$other = $wb.Worksheets.Item("Other")
$other.Move([Type]::Missing, $group2)

# --- Leave the remembered selection on "Other" (E36) --------------------------
$other = $wb.Worksheets.Item("Other")
[void]$other.Activate()
[void]$other.Range("E36").Select()

# --- Finish back on "Group 1" with the cursor at C28 --------------------------
$group1 = $wb.Worksheets.Item("Group 1")
[void]$group1.Activate()
[void]$group1.Range("C28").Select()
